$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Recolor row 5 to match the formatting of row 6 (cluster color), then clear the
# "OLED CS" (C5) and "can be grounded (P0_2=0)" (F5) text that no longer applies.
$src = $ws.Range("B6:F6")
$dst = $ws.Range("B5:F5")
$src.Copy()
$dst.PasteSpecial(-4122)
$ws.Range("C5").Value = ""
$ws.Range("F5").Value = ""

$excel.CutCopyMode = $false

$ws.Range("B11").Select() | Out-Null

Write-Host "done"
